{"js": "// Grammar fixes in the minutes document.\n// Each entry is an exact (old, new) full-sentence replacement, applied via\n// Body.search() + Range.insertText(..., \"Replace\") so the surrounding run\n// formatting (rPr) is preserved exactly like Word's Find & Replace would do.\nconst replacements = [\n  {\n    find: \"We have to find the solution on our own. To do that, we should conduct brainstorm or any other form of conversation.\",\n    replace: \"We have to find the solution on our own. To do that, we should conduct a brainstorming or starbursting session.\"\n  },\n  {\n    find: \"Another useful approach would be creating the graph representing position of every potential solution on two different axis: impact and difficulty. That would lead us in right direction when deciding on particular solution\",\n    replace: \"Another useful approach would be creating a graph representing position of every potential solution on two different axis: impact and difficulty. That would lead us in right direction when deciding on particular solution\"\n  },\n  {\n    find: \"We have watched the recording of ourselves presenting on mid-term presentation, and discussed about it.\",\n    replace: \"We have watched the recording of ourselves presenting on mid-term presentation, and reflected upon it.\"\n  },\n  {\n    find: \"spectators, and providing less distracting form of time measurement than an colorful animation in center of slides.\",\n    replace: \"spectators, and provide a less distracting form of time measurement than the colorful animation in center of slides.\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Grammar fixes in the minutes document.\n# Each entry is an exact (old, new) full-sentence replacement, applied via\n# Word's Find & Replace (Range.Find.Execute) idiom so the surrounding run\n# formatting is preserved just as it would be in a real Word macro.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($OldText, $NewText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n\n    $found = $rng.Find.Execute($OldText, $false, $false, $false, $false, $false, $true, 1, $false, $NewText, 1)\n\n    if (-not $found) {\n        throw \"Text not found: $OldText\"\n    }\n}\n\nReplace-ExactText \"We have to find the solution on our own. To do that, we should conduct brainstorm or any other form of conversation.\" \"We have to find the solution on our own. To do that, we should conduct a brainstorming or starbursting session.\"\n\nReplace-ExactText \"Another useful approach would be creating the graph representing position of every potential solution on two different axis: impact and difficulty. That would lead us in right direction when deciding on particular solution\" \"Another useful approach would be creating a graph representing position of every potential solution on two different axis: impact and difficulty. That would lead us in right direction when deciding on particular solution\"\n\nReplace-ExactText \"We have watched the recording of ourselves presenting on mid-term presentation, and discussed about it.\" \"We have watched the recording of ourselves presenting on mid-term presentation, and reflected upon it.\"\n\nReplace-ExactText \"spectators, and providing less distracting form of time measurement than an colorful animation in center of slides.\" \"spectators, and provide a less distracting form of time measurement than the colorful animation in center of slides.\"\n"}
